# Refactoring to fit personal dataset: rename header labels on row 1
# (columns B through H) and move the active selection to L7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "âge au dx"
$ws.Range("C1").Value = "Diag_PSA"
$ws.Range("D1").Value = "Gleason primaire Bx"
$ws.Range("E1").Value = "Gleason secondaire Bx"
$ws.Range("F1").Value = "Stade clinique"
$ws.Range("G1").Value = "NbCtePositive"
$ws.Range("H1").Value = "NbCteNegative"

$ws.Range("L7").Select()
